# Daily attendance processing - 2025-12-30 10:01:31
# Normalizes the "Recorded By" (column G) entries so that the "System" /
# "system" token that currently appears first in a two-part, comma
# separated list is moved to the end of the list (e.g.
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"). The one
# three-part case ("System, system, backup@backdoor.com") instead has its
# two "system" tokens swapped so the lowercase one leads, keeping the
# backdoor address last.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) {
        continue
    }

    $parts = $value -split ', '

    if ($parts.Count -eq 2 -and $parts[0] -eq 'System' -and $parts[1] -ne 'backup@backdoor.com') {
        $cell.Value = "$($parts[1]), $($parts[0])"
    }
    elseif ($parts.Count -eq 3 -and $parts[0] -eq 'System' -and $parts[1] -eq 'system' -and $parts[2] -eq 'backup@backdoor.com') {
        $cell.Value = "system, System, backup@backdoor.com"
    }
}
